$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: restructure the paragraphs around the "_GoBack" bookmark.
#
# Before: <p1 "Relating..."> <p2 empty> <p3 "Intro...">  ... <p9 bookmark+"How I want..."> ...
# After : <p1 "Relating..."> <p2 bookmark-only, empty>   <p3 "Intro...">  ... <p9 "How I want..." (no bookmark)> ...
# ---------------------------------------------------------------------

# 1a. Remove the bookmark from its old location (start of the "How I want
#     to do the research" paragraph).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 1b. Remove the old empty paragraph that sits right after the
#     "Relating Other Works" paragraph.
$oldEmptyPara = $d.Paragraphs.Item(2)
$oldEmptyPara.Range.Delete()

# 1c. Insert a fresh empty paragraph in the same spot (right after the
#     "Relating Other Works" paragraph).
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()

# 1d. Re-create the "_GoBack" bookmark so that both its start and end
#     markers live inside that new empty paragraph. A transient marker
#     character is used so the bookmark range is non-degenerate, then the
#     marker text is removed again, leaving only the bookmark behind.
$newEmptyPara = $d.Paragraphs.Item(2)
$newEmptyPara.Range.InsertAfter("X")
$markerRange = $d.Range($newEmptyPara.Range.Start, $newEmptyPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange = $d.Range($newEmptyPara.Range.Start, $newEmptyPara.Range.Start + 1)
$markerRange.Text = ""

# ---------------------------------------------------------------------
# Part 2: annotate the planning paragraph ("Intro 3 - 4 ..." / Background
# / Describe the approach / Related Works / Summary Conclusions) with the
# window-size notes, each as its own run right after the anchor text.
# ---------------------------------------------------------------------

function Insert-AfterText($searchText, $insertText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Collapse(0)
    $rng.InsertAfter($insertText)
    $newRng = $d.Range($rng.Start, $rng.Start + $insertText.Length)
    # Force the freshly inserted text into its own run (rather than being
    # silently merged back into the preceding run) by toggling a
    # character property on and back off.
    $newRng.Bold = 1
    $newRng.Bold = 0
}

Insert-AfterText "Intro 3 - 4 " " (3)"
Insert-AfterText "+ " " (9)"
Insert-AfterText "Describe the approach ~10" " (10)"
Insert-AfterText "Related Works 6-7" " (4)"
Insert-AfterText "Summary Conclusions 2" " (1)"
